# Append three new trade rows (4, 5, 6) to the IBB random-trade log sheet,
# mirroring the existing row 3 layout:
#   A=Principle  B=StartPrinciple  C=BuyPrice  D=SellPrice
#   E=IsShortSell  F=PriceChange%  G=Date  H=Profitable

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 9892.2900000000009
$ws.Range("B4").Value = 9946
$ws.Range("C4").Value = 297.77999999999997
$ws.Range("D4").Value = 296.18
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = -0.54
$ws.Range("G4").Value = 42606.585706018515
$ws.Range("H4").Value = $false

# Row 5
$ws.Range("A5").Value = 9838.8700000000008
$ws.Range("B5").Value = 9892.2900000000009
$ws.Range("C5").Value = 297.77999999999997
$ws.Range("D5").Value = 296.18
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = -0.54
$ws.Range("G5").Value = 42606.586863425924
$ws.Range("H5").Value = $false

# Row 6
$ws.Range("A6").Value = 9785.74
$ws.Range("B6").Value = 9838.8700000000008
$ws.Range("C6").Value = 297.77999999999997
$ws.Range("D6").Value = 296.18
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = -0.54
$ws.Range("G6").Value = 42606.58792824074
$ws.Range("H6").Value = $false
